$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("events")

$ws.Range("A6:L6").ClearFormats()

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Huánuco"
$ws.Range("C6").Value = "Lauricocha"
$ws.Range("D6").Value = "San Miguel de Cauri"
$ws.Range("E6").Value = 45987
$ws.Range("F6").Value = "Vientos fuertes"
$ws.Range("G6").Value = "En monitoreo"
$ws.Range("H6").Value = "En monitoreo"
$ws.Range("I6").Value = "En proceso"
$ws.Range("J6").Value = "Evaluación en curso"
$ws.Range("K6").Value = -10.298932000000001
$ws.Range("L6").Value = -76.637589000000006

$ws.Range("J9").Select()
